# Update the cryptos price list (Price / Volume(1h) refresh, plus a
# FraxShare / TrustWalletToken row swap) as produced by the scheduled
# GitHub Actions data refresh.
#
# Note: several "Price" values look numeric (e.g. "1.001", "0.8300") but
# must stay as literal text, matching how the source data was stored
# (leading zeros / trailing zeros preserved). A leading apostrophe forces
# Excel to keep them as text instead of silently converting to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.988.82"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.897.21"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'0.8300"
$ws.Range("E5").Value = "  +4.54%  "
$ws.Range("D6").Value = "'241.95"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D8").Value = "'0.3269"
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("D9").Value = "'26.40"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "'0.07016"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.08083"
$ws.Range("D12").Value = "'0.7622"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "1.895.05"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "'5.240"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'92.18"
$ws.Range("D16").Value = "29.989.57"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "'14.08"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'5.836"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "'243.22"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").Value = "'0.000007751"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "2.152.16"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'6.930"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'0.1726"
$ws.Range("E25").Value = "  +23.70%  "
$ws.Range("D26").Value = "'9.266"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'165.13"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").Value = "'18.90"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "'2.092"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").Value = "'1.359"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").Value = "'1.515"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").Value = "'0.05878"
$ws.Range("E32").Value = "  +7.97%  "
$ws.Range("D33").Value = "'4.274"
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").Value = "'4.065"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").Value = "'1.263"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'0.7291"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").Value = "'0.01918"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "'0.4428"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "'72.47"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.866"
$ws.Range("E42").Value = "  -5.04%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8548"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'1.897"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'101.73"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "'7.541"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "'9.779"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "'992.94"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").Value = "2.045.13"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "'1.519"
$ws.Range("E51").Value = "  +0.62%  "
